# Insert a new "Author" paragraph ("Ben Jarman") directly after the
# Subtitle paragraph ("Reflecting on ethics, consent, and reproducibility")
# and before the Date paragraph.

$d = $word.ActiveDocument

# Locate the Subtitle paragraph by its known text.
$subtitle = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Reflecting on ethics, consent, and reproducibility") {
        $subtitle = $p
        break
    }
}

# Insert a new paragraph right after the subtitle paragraph.
$subtitle.Range.InsertParagraphAfter()

# The newly created paragraph is now the one immediately following the
# subtitle paragraph; give it the "Author" style and its text.
$authorPara = $subtitle.Next()
$authorPara.Style = "Author"
$authorPara.Range.Text = "Ben Jarman"
